$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7483.85
$ws.Range("J17").Value = 7777.7896
$ws.Range("L17").Value = 23333.3688
$ws.Range("N17").Value = -23669.3688

$ws.Range("H106").Value = 9526752
$ws.Range("I106").Value = 37039096
$ws.Range("J106").Value = 3248.1538
$ws.Range("K106").Value = 37039096
$ws.Range("L106").Value = 3248.1538
$ws.Range("M106").Value = -37038465
$ws.Range("N106").Value = -4510.1538

$ws.Range("H111").Value = 2352.5454
$ws.Range("J111").Value = 898.5714
$ws.Range("L111").Value = 2695.7142
$ws.Range("N111").Value = -8829.7142

$ws.Range("H132").Value = 5474.5557
$ws.Range("I132").Value = 5736.1333
$ws.Range("J132").Value = 4166.6665
$ws.Range("K132").Value = 17208.3999
$ws.Range("L132").Value = 12499.9995
$ws.Range("M132").Value = -14678.3999
$ws.Range("N132").Value = -17559.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1416.6666
$ws.Range("I2").Value = 1425
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1425
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -1312
$ws.Range("N2").Value = -1626

$ws.Range("H32").Value = 5326.25
$ws.Range("I32").Value = 4570.476
$ws.Range("J32").Value = 9294.0625
$ws.Range("K32").Value = 4570.476
$ws.Range("L32").Value = 9294.0625
$ws.Range("M32").Value = -4283.476
$ws.Range("N32").Value = -9868.0625

$ws.Range("H61").Value = 1699.0227
$ws.Range("I61").Value = 1407.4
$ws.Range("J61").Value = 2833.111
$ws.Range("K61").Value = 1407.4
$ws.Range("L61").Value = 2833.111
$ws.Range("M61").Value = -1195.4
$ws.Range("N61").Value = -3257.111

$ws.Range("H74").Value = 30304552
$ws.Range("I74").Value = 43478924
$ws.Range("J74").Value = 3495.4
$ws.Range("K74").Value = 43478924
$ws.Range("L74").Value = 3495.4
$ws.Range("M74").Value = -43478050
$ws.Range("N74").Value = -5243.4

$ws.Range("H77").Value = 30304552
$ws.Range("I77").Value = 43478924
$ws.Range("J77").Value = 3495.4
$ws.Range("K77").Value = 217394620
$ws.Range("L77").Value = 17477
$ws.Range("M77").Value = -217390252
$ws.Range("N77").Value = -26213

$ws.Range("H116").Value = 1416.6666
$ws.Range("I116").Value = 1425
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 1425
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = 869
$ws.Range("N116").Value = -5988

$ws.Range("H132").Value = 12466.167
$ws.Range("I132").Value = 1844.7028
$ws.Range("J132").Value = 48192.91
$ws.Range("K132").Value = 5534.1084
$ws.Range("L132").Value = 144578.73
$ws.Range("M132").Value = -3004.1084
$ws.Range("N132").Value = -149638.73

$ws.Range("H136").Value = 1699.0227
$ws.Range("I136").Value = 1407.4
$ws.Range("J136").Value = 2833.111
$ws.Range("K136").Value = 4222.200000000001
$ws.Range("L136").Value = 8499.332999999999
$ws.Range("M136").Value = -1672.200000000001
$ws.Range("N136").Value = -13599.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1416.6666
$ws.Range("I3").Value = 1425
$ws.Range("J3").Value = 1400
$ws.Range("K3").Value = 1425
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = -1311
$ws.Range("N3").Value = -1628

$ws.Range("H20").Value = 1784.8846
$ws.Range("I20").Value = 1830.1818
$ws.Range("J20").Value = 1751.6666
$ws.Range("K20").Value = 1830.1818
$ws.Range("L20").Value = 1751.6666
$ws.Range("M20").Value = -1583.1818
$ws.Range("N20").Value = -2245.6666

$ws.Range("H107").Value = 1561.826
$ws.Range("I107").Value = 1559.0526
$ws.Range("J107").Value = 1575
$ws.Range("K107").Value = 1559.0526
$ws.Range("L107").Value = 1575
$ws.Range("M107").Value = 360.9474
$ws.Range("N107").Value = -5415

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3732.8364
$ws.Range("I31").Value = 1942.5652
$ws.Range("J31").Value = 5019.5938
$ws.Range("K31").Value = 1942.5652
$ws.Range("L31").Value = 5019.5938
$ws.Range("M31").Value = -1647.5652
$ws.Range("N31").Value = -5609.5938

$ws.Range("H34").Value = 3732.8364
$ws.Range("I34").Value = 1942.5652
$ws.Range("J34").Value = 5019.5938
$ws.Range("K34").Value = 1942.5652
$ws.Range("L34").Value = 5019.5938
$ws.Range("M34").Value = -1740.5652
$ws.Range("N34").Value = -5423.5938

$ws.Range("H43").Value = 32500
$ws.Range("J43").Value = 32500
$ws.Range("L43").Value = 32500
$ws.Range("N43").Value = -32868

$ws.Range("H86").Value = 59547.668
$ws.Range("I86").Value = 38585
$ws.Range("K86").Value = 38585
$ws.Range("M86").Value = -37462

$ws.Range("H89").Value = 59547.668
$ws.Range("I89").Value = 38585
$ws.Range("K89").Value = 192925
$ws.Range("M89").Value = -187309

$ws.Range("H101").Value = 32500
$ws.Range("J101").Value = 32500
$ws.Range("L101").Value = 32500
$ws.Range("N101").Value = -38990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7527.5
$ws.Range("J107").Value = 357.16666
$ws.Range("L107").Value = 1071.49998
$ws.Range("N107").Value = -4911.499980000001

$ws.Range("H113").Value = 713
$ws.Range("J113").Value = 716.4666999999999
$ws.Range("L113").Value = 2149.4001
$ws.Range("N113").Value = -6489.4001

$ws.Range("H122").Value = 826.1818
$ws.Range("J122").Value = 864
$ws.Range("L122").Value = 7776
$ws.Range("N122").Value = -12676

$ws.Range("H129").Value = 201177.6
$ws.Range("J129").Value = 236513.06
$ws.Range("L129").Value = 709539.1799999999
$ws.Range("N129").Value = -719539.1799999999

$ws.Range("H131").Value = 729.09
$ws.Range("J131").Value = 742.08246
$ws.Range("L131").Value = 2226.24738
$ws.Range("N131").Value = -12306.24738

$ws.Range("H132").Value = 864.6
$ws.Range("I132").Value = 746.25
$ws.Range("J132").Value = 999.8570999999999
$ws.Range("K132").Value = 6716.25
$ws.Range("L132").Value = 8998.713899999999
$ws.Range("M132").Value = -4186.25
$ws.Range("N132").Value = -14058.7139

$ws.Range("H136").Value = 3214.6155
$ws.Range("I136").Value = 1131.6666
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 3394.9998
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = 1705.0002
$ws.Range("N136").Value = -25200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6876.7144
$ws.Range("I102").Value = 10174.667
$ws.Range("J102").Value = 4403.25
$ws.Range("K102").Value = 10174.667
$ws.Range("L102").Value = 4403.25
$ws.Range("M102").Value = -8552.666999999999
$ws.Range("N102").Value = -7647.25

$ws.Range("H132").Value = 26762.043
$ws.Range("I132").Value = 5534.8887
$ws.Range("J132").Value = 103179.8
$ws.Range("K132").Value = 16604.6661
$ws.Range("L132").Value = 309539.4
$ws.Range("M132").Value = -14074.6661
$ws.Range("N132").Value = -314599.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 501515.5
$ws.Range("I35").Value = 501515.5
$ws.Range("K35").Value = 501515.5
$ws.Range("M35").Value = -501179.5

$ws.Range("H46").Value = 1551.9131
$ws.Range("I46").Value = 1313.7142
$ws.Range("J46").Value = 1922.4445
$ws.Range("K46").Value = 1313.7142
$ws.Range("L46").Value = 1922.4445
$ws.Range("M46").Value = -1125.7142
$ws.Range("N46").Value = -2298.4445

$ws.Range("H68").Value = 2628.6
$ws.Range("I68").Value = 2125
$ws.Range("J68").Value = 2964.3333
$ws.Range("K68").Value = 2125
$ws.Range("L68").Value = 2964.3333
$ws.Range("M68").Value = -1376
$ws.Range("N68").Value = -4462.3333

$ws.Range("H71").Value = 2628.6
$ws.Range("I71").Value = 2125
$ws.Range("J71").Value = 2964.3333
$ws.Range("K71").Value = 10625
$ws.Range("L71").Value = 14821.6665
$ws.Range("M71").Value = -6881
$ws.Range("N71").Value = -22309.6665

$ws.Range("H82").Value = 1207.3334
$ws.Range("I82").Value = 1269.091
$ws.Range("J82").Value = 1037.5
$ws.Range("K82").Value = 1269.091
$ws.Range("L82").Value = 1037.5
$ws.Range("M82").Value = -908.0909999999999
$ws.Range("N82").Value = -1759.5

$ws.Range("H85").Value = 1207.3334
$ws.Range("I85").Value = 1269.091
$ws.Range("J85").Value = 1037.5
$ws.Range("K85").Value = 1269.091
$ws.Range("L85").Value = 1037.5
$ws.Range("M85").Value = -21.09099999999989
$ws.Range("N85").Value = -3533.5

$ws.Range("H93").Value = 1494.36
$ws.Range("I93").Value = 1383.9048
$ws.Range("J93").Value = 2074.25
$ws.Range("K93").Value = 1383.9048
$ws.Range("L93").Value = 2074.25
$ws.Range("M93").Value = -135.9048
$ws.Range("N93").Value = -4570.25

$ws.Range("H100").Value = 2298.3333
$ws.Range("I100").Value = 1151
$ws.Range("J100").Value = 2527.8
$ws.Range("K100").Value = 1151
$ws.Range("L100").Value = 2527.8
$ws.Range("M100").Value = -610
$ws.Range("N100").Value = -3609.8

$ws.Range("H122").Value = 2455740.5
$ws.Range("I122").Value = 2805460.5
$ws.Range("K122").Value = 8416381.5
$ws.Range("M122").Value = -8413931.5

$ws.Range("H132").Value = 465647
$ws.Range("J132").Value = 4063.1667
$ws.Range("L132").Value = 12189.5001
$ws.Range("N132").Value = -17249.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1884.2727
$ws.Range("J81").Value = 3004.3333
$ws.Range("L81").Value = 6008.6666
$ws.Range("N81").Value = -8130.6666

$ws.Range("H84").Value = 1884.2727
$ws.Range("J84").Value = 3004.3333
$ws.Range("L84").Value = 30043.333
$ws.Range("N84").Value = -40651.333

$ws.Range("H132").Value = 959
$ws.Range("I132").Value = 876.8823
$ws.Range("J132").Value = 1041.1177
$ws.Range("K132").Value = 2630.6469
$ws.Range("L132").Value = 3123.3531
$ws.Range("M132").Value = -100.6468999999997
$ws.Range("N132").Value = -8183.3531

$ws.Range("H136").Value = 24580320
$ws.Range("I136").Value = 32259298
$ws.Range("J136").Value = 7590
$ws.Range("K136").Value = 96777894
$ws.Range("L136").Value = 22770
$ws.Range("M136").Value = -96775344
$ws.Range("N136").Value = -27870
